$wb = $excel.ActiveWorkbook

# Sheet "weibull"
$ws = $wb.Worksheets.Item("weibull")
$ws.Range("B2").Value = -2.39500202175584
$ws.Range("C2").Value = 0.10812372759168
$ws.Range("B3").Value = 0.103772517733394
$ws.Range("C3").Value = 0.0769739957034721

# Sheet "lognormal"
$ws = $wb.Worksheets.Item("lognormal")
$ws.Range("B2").Value = 2.03645638169679
$ws.Range("C2").Value = 0.152717764270313
$ws.Range("B3").Value = -1.12190166059483
$ws.Range("C3").Value = 0.0867906653867666

# Sheet "llogis"
$ws = $wb.Worksheets.Item("llogis")
$ws.Range("B2").Value = -1.777951348557
$ws.Range("C2").Value = 0.0825379900743899
$ws.Range("B3").Value = 0.679424176932772
$ws.Range("C3").Value = 0.0853173482110215

# Sheet "gompertz"
$ws = $wb.Worksheets.Item("gompertz")
$ws.Range("B2").Value = -2.08140174659871
$ws.Range("C2").Value = 0.101464536652
$ws.Range("B3").Value = -0.0180406304756775
$ws.Range("C3").Value = 0.012319060510442

# Sheet "weibull cov"
$ws = $wb.Worksheets.Item("weibull cov")
$ws.Range("A2").Value = 0.0116907404683198
$ws.Range("B2").Value = -0.00402251489432989
$ws.Range("A3").Value = -0.00402251489432989
$ws.Range("B3").Value = 0.00592499601455814

# Sheet "lognormal cov"
$ws = $wb.Worksheets.Item("lognormal cov")
$ws.Range("A2").Value = 0.0233227155237229
$ws.Range("B2").Value = -0.0107506639965146
$ws.Range("A3").Value = -0.0107506639965146
$ws.Range("B3").Value = 0.00753261959827769

# Sheet "llogis cov"
$ws = $wb.Worksheets.Item("llogis cov")
$ws.Range("A2").Value = 0.00681251980552008
$ws.Range("B2").Value = 0.00220984086502586
$ws.Range("A3").Value = 0.00220984086502586
$ws.Range("B3").Value = 0.00727904990576069

# Sheet "gompertz cov"
$ws = $wb.Worksheets.Item("gompertz cov")
$ws.Range("A2").Value = 0.0102950521980051
$ws.Range("B2").Value = -0.000601428257233394
$ws.Range("A3").Value = -0.000601428257233394
$ws.Range("B3").Value = 0.000151759251859932
